$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "93-19=",
    "14+12=",
    "8+8=",
    "72+16=",
    "68-65=",
    "21+11=",
    "21+24=",
    "61-45=",
    "65-5=",
    "23+27=",
    "90-1=",
    "97-89=",
    "25+25=",
    "46-42=",
    "48+42=",
    "73-25=",
    "16-12=",
    "82-21=",
    "71-52=",
    "84+12=",
    "59-24=",
    "58-3=",
    "40+2=",
    "29-15=",
    "45+8=",
    "97-57=",
    "21+72=",
    "51-8=",
    "37+33=",
    "18+40=",
    "57-5=",
    "43+56=",
    "17-1=",
    "26+8=",
    "3+94=",
    "60+5=",
    "96-94=",
    "36+28=",
    "76-43=",
    "84-35=",
    "59+1=",
    "63+23=",
    "83+4=",
    "65-18=",
    "39+54=",
    "27+66=",
    "38+50=",
    "74+13=",
    "79+11=",
    "45+27=",
    "42-39=",
    "48-21=",
    "21+38=",
    "38+0=",
    "46+36=",
    "13+37=",
    "94-92=",
    "79-28=",
    "19-6=",
    "94-47=",
    "21+65=",
    "5+18=",
    "52+38=",
    "48-37=",
    "8+86=",
    "74-51=",
    "80+2=",
    "31+61=",
    "40-15=",
    "65+23=",
    "1+27=",
    "69-21=",
    "21+15=",
    "37-9=",
    "91-34=",
    "96-64=",
    "25+5=",
    "86-1=",
    "38+30=",
    "35+26=",
    "24-1=",
    "24+69=",
    "6+0=",
    "52-23=",
    "44+5=",
    "1+83=",
    "78-10=",
    "99-14=",
    "0+98=",
    "77-51=",
    "7+75=",
    "73+7=",
    "81-15=",
    "27-26=",
    "25+45=",
    "20+74=",
    "31+24=",
    "43+24=",
    "91-70=",
    "45+9="
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Host "Updated $idx cells"
